$wb = $excel.ActiveWorkbook

# The existing "Czech" sheet is the template for every country tab in this
# workbook (same layout/styles, just different market name + NGC code), so
# build the new "Swiss" tab by copying it and patching the two data cells.
$czech = $wb.Worksheets.Item("Czech")
$czech.Select()
$czech.Copy($null, $czech) | Out-Null

$swiss = $wb.Worksheets.Item($czech.Index + 1)
$swiss.Name = "Swiss"

$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2647/T2657"

# Czech itself is left with a "select all columns" state (as seen through the
# UI when flipping to the newly duplicated tab), before landing on the new tab.
$czech.Select()
$czech.Cells.EntireColumn.Select() | Out-Null

$swiss.Select()
$swiss.Range("B7").Select() | Out-Null
